$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Tnfrsf14"
$ws.Range("C2").Value = "Cd160"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 6.810423666666666
$ws.Range("H2").Value = 20.431271
$ws.Range("I2").Value = 0.4869218264300073
$ws.Range("J2").Value = 0.4869218264300074
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.106442666666667
$ws.Range("N2").Value = 6.319328000000001
$ws.Range("O2").Value = 0.6232646226081526
$ws.Range("P2").Value = 0.6232646226081526
$ws.Range("Q2").Value = 14.34576698954311
$ws.Range("R2").Value = 129.111902905888
$ws.Range("S2").Value = 0.3034811483895709
$ws.Range("T2").Value = 0.303481148389571

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Tnfrsf14"
$ws.Range("C3").Value = "Cd160"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 6.810423666666666
$ws.Range("H3").Value = 20.431271
$ws.Range("I3").Value = 0.4869218264300073
$ws.Range("J3").Value = 0.4869218264300074
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.6936943333333333
$ws.Range("N3").Value = 2.081083
$ws.Range("O3").Value = 0.2052536932109303
$ws.Range("P3").Value = 0.2052536932109303
$ws.Range("Q3").Value = 4.724352305165888
$ws.Range("R3").Value = 42.519170746493
$ws.Range("S3").Value = 0.09994250317977059
$ws.Range("T3").Value = 0.0999425031797706

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Tnfrsf14"
$ws.Range("C4").Value = "Cd160"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 6.810423666666666
$ws.Range("H4").Value = 20.431271
$ws.Range("I4").Value = 0.4869218264300073
$ws.Range("J4").Value = 0.4869218264300074
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.5795553333333333
$ws.Range("N4").Value = 1.738666
$ws.Range("O4").Value = 0.1714816841809171
$ws.Range("P4").Value = 0.1714816841809171
$ws.Range("Q4").Value = 3.947017358276222
$ws.Range("R4").Value = 35.523156224486
$ws.Range("S4").Value = 0.08349817486066582
$ws.Range("T4").Value = 0.08349817486066584

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Tnfrsf14"
$ws.Range("C5").Value = "Cd160"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.383140666666667
$ws.Range("H5").Value = 16.149422
$ws.Range("I5").Value = 0.3848760097220062
$ws.Range("J5").Value = 0.3848760097220062
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.106442666666667
$ws.Range("N5").Value = 6.319328000000001
$ws.Range("O5").Value = 0.6232646226081526
$ws.Range("P5").Value = 0.6232646226081526
$ws.Range("Q5").Value = 11.33927718093511
$ws.Range("R5").Value = 102.053494628416
$ws.Range("S5").Value = 0.2398796009503179
$ws.Range("T5").Value = 0.2398796009503179

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Tnfrsf14"
$ws.Range("C6").Value = "Cd160"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.383140666666667
$ws.Range("H6").Value = 16.149422
$ws.Range("I6").Value = 0.3848760097220062
$ws.Range("J6").Value = 0.3848760097220062
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.6936943333333333
$ws.Range("N6").Value = 2.081083
$ws.Range("O6").Value = 0.2052536932109303
$ws.Range("P6").Value = 0.2052536932109303
$ws.Range("Q6").Value = 3.734254176002889
$ws.Range("R6").Value = 33.608287584026
$ws.Range("S6").Value = 0.0789972224237277
$ws.Range("T6").Value = 0.0789972224237277

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Tnfrsf14"
$ws.Range("C7").Value = "Cd160"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.383140666666667
$ws.Range("H7").Value = 16.149422
$ws.Range("I7").Value = 0.3848760097220062
$ws.Range("J7").Value = 0.3848760097220062
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.5795553333333333
$ws.Range("N7").Value = 1.738666
$ws.Range("O7").Value = 0.1714816841809171
$ws.Range("P7").Value = 0.1714816841809171
$ws.Range("Q7").Value = 3.119827883450222
$ws.Range("R7").Value = 28.078450951052
$ws.Range("S7").Value = 0.06599918634796063
$ws.Range("T7").Value = 0.06599918634796063

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Tnfrsf14"
$ws.Range("C8").Value = "Cd160"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.793123666666667
$ws.Range("H8").Value = 5.379371
$ws.Range("I8").Value = 0.1282021638479865
$ws.Range("J8").Value = 0.1282021638479865
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.106442666666667
$ws.Range("N8").Value = 6.319328000000001
$ws.Range("O8").Value = 0.6232646226081526
$ws.Range("P8").Value = 0.6232646226081526
$ws.Range("Q8").Value = 3.777112198076444
$ws.Range("R8").Value = 33.994009782688
$ws.Range("S8").Value = 0.07990387326826386
$ws.Range("T8").Value = 0.07990387326826387

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Tnfrsf14"
$ws.Range("C9").Value = "Cd160"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.793123666666667
$ws.Range("H9").Value = 5.379371
$ws.Range("I9").Value = 0.1282021638479865
$ws.Range("J9").Value = 0.1282021638479865
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.6936943333333333
$ws.Range("N9").Value = 2.081083
$ws.Range("O9").Value = 0.2052536932109303
$ws.Range("P9").Value = 0.2052536932109303
$ws.Range("Q9").Value = 1.243879726532555
$ws.Range("R9").Value = 11.194917538793
$ws.Range("S9").Value = 0.02631396760743205
$ws.Range("T9").Value = 0.02631396760743205

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Tnfrsf14"
$ws.Range("C10").Value = "Cd160"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.793123666666667
$ws.Range("H10").Value = 5.379371
$ws.Range("I10").Value = 0.1282021638479865
$ws.Range("J10").Value = 0.1282021638479865
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.5795553333333333
$ws.Range("N10").Value = 1.738666
$ws.Range("O10").Value = 0.1714816841809171
$ws.Range("P10").Value = 0.1714816841809171
$ws.Range("Q10").Value = 1.039214384342889
$ws.Range("R10").Value = 9.352929459086
$ws.Range("S10").Value = 0.0219843229722906
$ws.Range("T10").Value = 0.02198432297229061

